# Doing full stat run
# Updates the "Runs" sheet H-column benchmark timings (re-measured numbers
# from a full stat run), tags row 98 with the "xxxx" marker in column I,
# and switches the active/selected sheet from "Simple Test" to "Runs".

$wb = $excel.ActiveWorkbook
$wsRuns = $wb.Worksheets.Item("Runs")
$wsSimple = $wb.Worksheets.Item("Simple Test")

$wsRuns.Cells.Item(2, 8).Value = 17
$wsRuns.Cells.Item(5, 8).Value = 27
$wsRuns.Cells.Item(6, 8).Value = 46
$wsRuns.Cells.Item(7, 8).Value = 49
$wsRuns.Cells.Item(8, 8).Value = 32
$wsRuns.Cells.Item(9, 8).Value = 51
$wsRuns.Cells.Item(10, 8).Value = 55
$wsRuns.Cells.Item(11, 8).Value = 29
$wsRuns.Cells.Item(12, 8).Value = 48
$wsRuns.Cells.Item(13, 8).Value = 51
$wsRuns.Cells.Item(14, 8).Value = 18
$wsRuns.Cells.Item(17, 8).Value = 28
$wsRuns.Cells.Item(18, 8).Value = 47
$wsRuns.Cells.Item(19, 8).Value = 52
$wsRuns.Cells.Item(20, 8).Value = 33
$wsRuns.Cells.Item(21, 8).Value = 54
$wsRuns.Cells.Item(22, 8).Value = 57
$wsRuns.Cells.Item(23, 8).Value = 30
$wsRuns.Cells.Item(25, 8).Value = 52
$wsRuns.Cells.Item(27, 8).Value = 47
$wsRuns.Cells.Item(28, 8).Value = 51
$wsRuns.Cells.Item(29, 8).Value = 29
$wsRuns.Cells.Item(30, 8).Value = 70
$wsRuns.Cells.Item(31, 8).Value = 74
$wsRuns.Cells.Item(32, 8).Value = 33
$wsRuns.Cells.Item(33, 8).Value = 74
$wsRuns.Cells.Item(34, 8).Value = 78
$wsRuns.Cells.Item(35, 8).Value = 34
$wsRuns.Cells.Item(36, 8).Value = 64
$wsRuns.Cells.Item(37, 8).Value = 68
$wsRuns.Cells.Item(38, 8).Value = 27
$wsRuns.Cells.Item(39, 8).Value = 46
$wsRuns.Cells.Item(40, 8).Value = 48
$wsRuns.Cells.Item(41, 8).Value = 27
$wsRuns.Cells.Item(63, 8).Value = 34
$wsRuns.Cells.Item(64, 8).Value = 35
$wsRuns.Cells.Item(66, 8).Value = 53
$wsRuns.Cells.Item(67, 8).Value = 54
$wsRuns.Cells.Item(69, 8).Value = 60
$wsRuns.Cells.Item(70, 8).Value = 60
$wsRuns.Cells.Item(72, 8).Value = 59
$wsRuns.Cells.Item(73, 8).Value = 60
$wsRuns.Cells.Item(74, 8).Value = 12
$wsRuns.Cells.Item(75, 8).Value = 36
$wsRuns.Cells.Item(76, 8).Value = 36
$wsRuns.Cells.Item(78, 8).Value = 56
$wsRuns.Cells.Item(79, 8).Value = 56
$wsRuns.Cells.Item(81, 8).Value = 63
$wsRuns.Cells.Item(82, 8).Value = 63
$wsRuns.Cells.Item(83, 8).Value = 24
$wsRuns.Cells.Item(84, 8).Value = 61
$wsRuns.Cells.Item(85, 8).Value = 60
$wsRuns.Cells.Item(89, 8).Value = 21
$wsRuns.Cells.Item(90, 8).Value = 78
$wsRuns.Cells.Item(92, 8).Value = 26
$wsRuns.Cells.Item(93, 8).Value = 86
$wsRuns.Cells.Item(94, 8).Value = 89
$wsRuns.Cells.Item(95, 8).Value = 25
$wsRuns.Cells.Item(96, 8).Value = 85
$wsRuns.Cells.Item(97, 8).Value = 89

# New "xxxx" tag added next to row 98 (column I)
$wsRuns.Cells.Item(98, 9).Value = "xxxx"

# Simple Test" no longer scrolled/selected at J22 -> now parked at A22
$wsSimple.Activate()
$wsSimple.Range("A22").Select()

# Make "Runs" the active sheet/tab, with H2 selected (matches activeTab=1 in workbook.xml)
$wsRuns.Activate()
$wsRuns.Range("H2").Select()

